# Insert three new "ListParagraph" bullet items right after the paragraph
# that reads "Diff b/w interpolation and property binding" and right before
# the (empty) paragraph that terminates the document body.
#
#   * "If-Else"                       (ilvl 0, numId 1)
#   * "ngIf directive" (with proofErr markers around "ngIf")  (ilvl 1, numId 1)
#   * "ng-template "                  (ilvl 1, numId 1, trailing space)

$d = $word.ActiveDocument

# Locate the anchor paragraph by its exact (trimmed) text.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Diff b/w interpolation and property binding") {
        $anchor = $p
    }
}

if ($anchor -eq $null) {
    Write-Host "Anchor paragraph not found!"
} else {
    # The paragraph's Range includes its trailing paragraph-mark character,
    # so End-1 is the character position immediately BEFORE that mark --
    # i.e. exactly where new paragraphs should be spliced in without
    # consuming (and thus replacing) the following paragraph.
    $insertAt = $anchor.Range.End - 1
    $rng = $d.Range($insertAt, $insertAt)

    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>If-Else</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ngIf</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> directive</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">ng-template </w:t></w:r></w:p>'

    $rng.InsertXML($xml)
}
